$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.03
$ws.Range("E6").Value = 13.045
$ws.Range("C7").Value = -13.293
$ws.Range("E7").Value = 13.045
$ws.Range("A8").Value = -21.153
$ws.Range("E8").Value = 12.919
$ws.Range("E9").Value = 12.3
$ws.Range("A10").Value = -20.712
$ws.Range("E10").Value = 12.518
$ws.Range("A12").Value = -21.649
$ws.Range("E12").Value = 13.055
$ws.Range("B13").Value = 6.595000000000001
$ws.Range("A18").Value = -21.649
$ws.Range("C20").Value = -12.976
$ws.Range("A25").Value = -21.754
